$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new test case row is being inserted at row 23 ("List of products added into
# shopping cart" / TC4 / Shopping cart section), which pushes every following
# row down by one (old row 23 -> new row 24, ..., old row 32 -> new row 33).
#
# We replicate that row-insert + shift manually (rather than via Rows.Insert())
# because a plain Insert() next to this particular data causes the COM runtime
# to fabricate a brand-new (unused) cell style in styles.xml. Doing the shift as
# explicit copy/paste of existing rows re-uses the existing style indices, which
# matches the original file's style table.

# Step 1: materialize a new blank row 33 at the bottom of the table by inserting
# next to two uniformly-styled blank rows (30-32) -- this does not introduce a
# new style -- then copy row 32's formatting onto it.
$ws.Rows.Item(33).Insert()
$ws.Range("A33:F33").ClearContents()
$ws.Range("A32:F32").Copy($ws.Range("A33:F33"))

# Step 2: shift the contents of rows 23-32 down into rows 24-33, working from
# the bottom up so that source rows are not overwritten before they are read.
for ($r = 32; $r -ge 23; $r--) {
    $dstRow = $r + 1
    $ws.Range("A" + $dstRow + ":F" + $dstRow).ClearContents()
    $ws.Range("A" + $r + ":F" + $r).Copy($ws.Range("A" + $dstRow + ":F" + $dstRow))
}

# Step 3: populate the now-empty row 23 with the new test case, re-using row
# 22's formatting (border everywhere, wrapped text in column C) which matches
# the desired style pattern for the new row.
$ws.Range("A23:F23").ClearContents()
$ws.Range("A22:F22").Copy($ws.Range("A23:F23"))

$ws.Range("A23").Value = "Shopping cart"
$ws.Range("B23").Value = "TC4"
$ws.Range("C23").Value = "List of products added into shopping cart"
$ws.Range("D23").Value = "Medium"
$ws.Range("E23").Value = "Sanity"
$ws.Range("F23").Value = ""

$excel.CutCopyMode = 0

# Update the active selection to match the edited area.
$ws.Range("E24").Select() | Out-Null
